# Applies the "updated table so it's better" changes to the workbook.

$wb = $excel.ActiveWorkbook

# 1. Rename "Backup" sheet to "Backups"
$backupSheet = $wb.Worksheets.Item("Backup")
$backupSheet.Name = "Backups"

# 2. Drive sheet: fix column name "Drive_Name" -> "Name", and update selection
$driveSheet = $wb.Worksheets.Item("Drive")
$driveSheet.Range("B3").Value = "Name"
$driveSheet.Select()
$driveSheet.Range("B4").Select()

# 3. User Properties sheet: remove the "Run on startup" row (row 4) and update selection
$userPropsSheet = $wb.Worksheets.Item("User Properties")
$userPropsSheet.Range("A4:E4").ClearContents()
$userPropsSheet.Select()
$userPropsSheet.Range("C17:C18").Select()
